$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.518.73'
$ws.Range("E2").Value = '  +5.82%  '
$ws.Range("D3").Value = '3.457.80'
$ws.Range("E3").Value = '  +4.15%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '414.66'
$ws.Range("E5").Value = '  +0.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.97'
$ws.Range("E6").Value = '  +17.58%  '
$ws.Range("D7").Value = '3.450.50'
$ws.Range("E7").Value = '  +3.98%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.592'
$ws.Range("E8").Value = '  +1.29%  '
$ws.Range("E9").Value = '  +0.25%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.689'
$ws.Range("E10").Value = '  +8.81%  '
$ws.Range("E11").Value = '  +29.65%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.69'
$ws.Range("E12").Value = '  +9.93%  '
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("D14").Value = '4.009.07'
$ws.Range("E14").Value = '  +4.51%  '
$ws.Range("E15").Value = '  +4.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.28'
$ws.Range("E16").Value = '  +4.61%  '
$ws.Range("D17").Value = '3.456.12'
$ws.Range("E17").Value = '  +4.15%  '
$ws.Range("D18").Value = '62.469.91'
$ws.Range("E18").Value = '  +6.16%  '
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.95'
$ws.Range("E20").Value = '  +1.97%  '
$ws.Range("E21").Value = '  +26.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.35'
$ws.Range("E22").Value = '  +1.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.25'
$ws.Range("E23").Value = '  +2.18%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '81.83'
$ws.Range("E24").Value = '  +9.16%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '313.52'
$ws.Range("E25").Value = '  +3.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.21'
$ws.Range("E26").Value = '  -0.77%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '30.29'
$ws.Range("E27").Value = '  +6.13%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.14'
$ws.Range("E28").Value = '  +3.66%  '
$ws.Range("B29").Value = 'Hedera'
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.123'
$ws.Range("E29").Value = '  +9.80%  '
$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.78'
$ws.Range("E30").Value = '  +6.37%  '
$ws.Range("E31").Value = '  +3.99%  '
$ws.Range("E32").Value = '  -1.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '45.08'
$ws.Range("E33").Value = '  +11.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.96'
$ws.Range("E34").Value = '  +4.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.62'
$ws.Range("E35").Value = '  +22.99%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0496'
$ws.Range("E37").Value = '  -6.19%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '52.75'
$ws.Range("E38").Value = '  +1.45%  '
$ws.Range("E39").Value = '  +3.28%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("E41").Value = '  -6.49%  '
$ws.Range("E42").Value = '  +6.03%  '
$ws.Range("E43").Value = '  +2.86%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '136.73'
$ws.Range("E44").Value = '  -0.62%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '17.76'
$ws.Range("E45").Value = '  +6.52%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.01'
$ws.Range("E46").Value = '  +2.23%  '
$ws.Range("B47").Value = 'TheGraph'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.289'
$ws.Range("E47").Value = '  +4.12%  '
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.49'
$ws.Range("E49").Value = '  +0.54%  '
$ws.Range("D50").Value = '2.249.12'
$ws.Range("E50").Value = '  +2.74%  '
$ws.Range("D51").Value = '3.810.13'
$ws.Range("E51").Value = '  +4.48%  '
